# Update "Paises" (countries) and "provincias Spain" data sheet
# as part of the 22 May 2020 07:05 refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last refreshed" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 07:05"

# 2) Ghana overtakes Oman in total cases -> swap their rows (64 & 65),
#    Ghana gets fresh numbers, Oman keeps its previous figures.
$ws.Range("A64").Value = "Ghana"
$ws.Range("B64").Value = 6486
$ws.Range("C64").Value = 217
$ws.Range("D64").Value = 1951
$ws.Range("E64").Value = 4504
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 31

$ws.Range("A65").Value = "Oman"
$ws.Range("B65").Value = 6370
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 1821
$ws.Range("E65").Value = 4518
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 31

# 3) Hungria (row 72) refreshed figures
$ws.Range("B72").Value = 3678
$ws.Range("C72").Value = 37
$ws.Range("D72").Value = 1587
$ws.Range("E72").Value = 1615
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 476

# 4) Tailandia (row 76) refreshed figures
$ws.Range("D76").Value = 2910
$ws.Range("E76").Value = 71

# 5) Uzbekistan (row 77) refreshed figures
$ws.Range("B77").Value = 2967
$ws.Range("C77").Value = 3
$ws.Range("E77").Value = 547
